$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- Sheet: About ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Linhuan Coal Mine, China, M1319, version ' + "'" + $newVersion + "'" + ". (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet: Boundaries and methane sources ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Find the last used row by scanning column A (mine name) from the top.
$lastRow = 1
while ($wsData.Cells.Item($lastRow + 1, 1).Value2 -ne $null) {
    $lastRow = $lastRow + 1
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # Column S = build_version
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
